$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Remove the worker row for period 2507 / doc 1037571366 / OSCAR DANIEL SOLAR MUÑOZ
# (whole row 17 is deleted; rows below shift up automatically).
$ws.Range("B17:J17").EntireRow.Delete()

# Update the account-statement totals for the remaining worker (LISNEYS PATRICIA BLANQUICET ROMERO):
# Valor Mora total goes from 113880 to 56940 (only one period of mora left)
$ws.Range("E11").Value = 56940

# Cant. Periodos goes from 2 to 1
$ws.Range("C13").Value = 1

# The remaining mora period advances from 2507 to 2508
$ws.Range("E16").Value = "2508"

Write-Host "edit applied"
